# Generate Report for Handoff
# Adds a new file (b4cd2c1f-232f-4643-b458-edbe832cb248.md) as a new row
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (table3 / Overview)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item("Overview")
$rowOverview = $loOverview.ListRows.Add()
$rOverview = $rowOverview.Range.Row

$wsOverview.Range("A" + $rOverview).Value = "b4cd2c1f-232f-4643-b458-edbe832cb248.md"
$wsOverview.Range("B" + $rOverview).Value = "e2e\b4cd2c1f-232f-4643-b458-edbe832cb248.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B" + $rOverview), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6827eed636b21d7975f2ae4177dc3bc2c5f1985f/e2e/b4cd2c1f-232f-4643-b458-edbe832cb248.md", "", "", "e2e\b4cd2c1f-232f-4643-b458-edbe832cb248.md")
$wsOverview.Range("C" + $rOverview).Value = ".md"
$wsOverview.Range("E" + $rOverview).Value = "Ready for handoff"
$wsOverview.Range("F" + $rOverview).Value = "Ready for handoff"
$wsOverview.Range("G" + $rOverview).Value = "2016-11-14 17:51:46"

# ---------------------------------------------------------------------
# Sheet "zh-cn" (table1 / zh_cn)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item("zh_cn")
$rowZhCn = $loZhCn.ListRows.Add()
$rZhCn = $rowZhCn.Range.Row

$wsZhCn.Range("A" + $rZhCn).Value = "b4cd2c1f-232f-4643-b458-edbe832cb248.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A" + $rZhCn), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6827eed636b21d7975f2ae4177dc3bc2c5f1985f/e2e/b4cd2c1f-232f-4643-b458-edbe832cb248.md", "", "", "b4cd2c1f-232f-4643-b458-edbe832cb248.md")
$wsZhCn.Range("B" + $rZhCn).Value = ".md"
$wsZhCn.Range("C" + $rZhCn).Value = "Ready for handoff"
$wsZhCn.Range("D" + $rZhCn).Value = "e2e"
$wsZhCn.Range("E" + $rZhCn).Value = "ht"
$wsZhCn.Range("F" + $rZhCn).Value = "False"
$wsZhCn.Range("G" + $rZhCn).Value = "b4cd2c1f-232f-4643-b458-edbe832cb248.6827eed636b21d7975f2ae4177dc3bc2c5f1985f.zh-cn.xlf"
$wsZhCn.Range("H" + $rZhCn).Value = "2016-11-14 17:51:30"
$wsZhCn.Range("K" + $rZhCn).Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M" + $rZhCn).Value = "True"
$wsZhCn.Range("O" + $rZhCn).Value = "False"

# ---------------------------------------------------------------------
# Sheet "de-de" (table2 / de_de)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item("de_de")
$rowDeDe = $loDeDe.ListRows.Add()
$rDeDe = $rowDeDe.Range.Row

$wsDeDe.Range("A" + $rDeDe).Value = "b4cd2c1f-232f-4643-b458-edbe832cb248.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A" + $rDeDe), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6827eed636b21d7975f2ae4177dc3bc2c5f1985f/e2e/b4cd2c1f-232f-4643-b458-edbe832cb248.md", "", "", "b4cd2c1f-232f-4643-b458-edbe832cb248.md")
$wsDeDe.Range("B" + $rDeDe).Value = ".md"
$wsDeDe.Range("C" + $rDeDe).Value = "Ready for handoff"
$wsDeDe.Range("D" + $rDeDe).Value = "e2e"
$wsDeDe.Range("E" + $rDeDe).Value = "ht"
$wsDeDe.Range("F" + $rDeDe).Value = "False"
$wsDeDe.Range("G" + $rDeDe).Value = "b4cd2c1f-232f-4643-b458-edbe832cb248.6827eed636b21d7975f2ae4177dc3bc2c5f1985f.de-de.xlf"
$wsDeDe.Range("H" + $rDeDe).Value = "2016-11-14 17:51:46"
$wsDeDe.Range("K" + $rDeDe).Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M" + $rDeDe).Value = "True"
$wsDeDe.Range("O" + $rDeDe).Value = "False"
